$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The hyperlink cell A6 used to display/store the "TUNTIKIRJAUKSET" wiki URL;
# it now points at the "Resurssit-ja-tyotunnit" wiki page instead. Updating the
# cell text (the hyperlink's display text) also naturally reorders the shared
# string table the same way the real author's save did (A1's unrelated text
# "Tunnit yhteensa 2.1-8.1" ends up sharing the slot vacated by the old URL).
$ws.Range("A6").Value2 = "https://github.com/DigiaMinions/Project/wiki/Resurssit-ja-ty%C3%B6tunnit"

# Column widths were nudged (Excel re-ran its "best fit" sizing after the new,
# longer hyperlink text was entered). The host only lets us drive column width
# through the character-based ColumnWidth property, which snaps to the
# nearest 1/6-character increment -- so we pick the input that lands on the
# closest achievable increment to the recorded target width for each column.
$ws.Columns.Item(1).ColumnWidth = 69.6666666666667
$ws.Columns.Item(2).ColumnWidth = 10.8333333333333
$ws.Columns.Item(5).ColumnWidth = 23.5
$ws.Columns.Item(6).ColumnWidth = 27.6666666666667

# The active selection moved from A12 to C16.
$ws.Range("C16").Select()
